$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
